# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets, matching the latest generated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 95
    4  = 1544
    5  = 589
    7  = 11227
    10 = 142
    11 = 334
    12 = 1078
    13 = 771
    14 = 12281
    15 = 12891
    22 = 67
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
